$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" date column (C) for rows 2-10 from 2023-10-22 to 2023-10-25
for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 3).Value = Get-Date -Year 2023 -Month 10 -Day 25 -Hour 0 -Minute 0 -Second 0
}
